# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") on Sheet1 with freshly recomputed strikeout
# values for each saved game row (rows 2-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 0
    9  = 0
    10 = 3
    11 = 1
    12 = 3
    13 = 1
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 2
    24 = 1
    25 = 3
    26 = 2
    27 = 1
    28 = 0
    29 = 1
    30 = 1
    31 = 2
    32 = 2
    33 = 0
    34 = 2
    35 = 1
    36 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
